$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.113.82'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.802.11'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.97'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5086'
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3842'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07697'
$ws.Range('E9').Value = '  -3.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.100'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.67'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.339'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.003'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.33'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.801.98'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.277'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.11'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001069'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06564'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.23'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.959'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.136.71'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.05'
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.248'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.56'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.426'
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.012.77'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.24'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.74'
$ws.Range('E30').Value = '  +3.84%  '
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.045'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.650'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.542'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07001'
$ws.Range('E35').Value = '  -3.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.981'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02346'
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2168'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.011'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.43'
$ws.Range('E40').Value = '  -5.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6117'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.150'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.23'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.300'
$ws.Range('E45').Value = '  -5.17%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5910'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.721'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.09'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.184'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.915'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06737'
$ws.Range('E51').Value = '  -1.26%  '
